$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.028.73'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").Value = '1.597.90'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.64'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3776'
$ws.Range("E7").Value = '  -1.46%  '
$ws.Range("E8").Value = '  -2.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.94'
$ws.Range("E9").Value = '  +0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.249'
$ws.Range("E10").Value = '  -3.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08124'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.31'
$ws.Range("E13").Value = '  -3.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.576'
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.340'
$ws.Range("E15").Value = '  -3.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001242'
$ws.Range("E16").Value = '  -2.80%  '
$ws.Range("D17").Value = '1.597.19'
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06811'
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("E20").Value = '  -4.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.509'
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.5577'
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.01'
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").Value = '23.032.94'
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.361'
$ws.Range("E26").Value = '  -1.95%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.791'
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.03'
$ws.Range("E28").Value = '  -2.25%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '149.78'
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.231'
$ws.Range("E30").Value = '  -3.24%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.09'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.386'
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.799'
$ws.Range("E33").Value = '  -13.52%  '
$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D34").Value = '1.772.19'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9577'
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07558'
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.24'
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.193'
$ws.Range("E38").Value = '  -4.32%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02689'
$ws.Range("E39").Value = '  -5.37%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2518'
$ws.Range("E40").Value = '  -3.28%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08815'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.366'
$ws.Range("E42").Value = '  -2.22%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7004'
$ws.Range("E43").Value = '  -4.33%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.35'
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.09'
$ws.Range("E45").Value = '  -5.55%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6604'
$ws.Range("E46").Value = '  -2.09%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.991'
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.274'
$ws.Range("E49").Value = '  -3.63%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.08'
$ws.Range("E50").Value = '  -0.92%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07898'
$ws.Range("E51").Value = '  -2.38%  '
